# Append: 2025-11-05 06:27 JST
# Update the "ランサーズ" (Lancers) listing sheet with a fresh scrape:
#  - refresh the "取得日時" timestamp on every kept row
#  - rows 5-8 get new listings (shifted in from the newer scrape)
#  - rows 9-19 (the old, now-stale listings) are dropped entirely
#  - column D / H get a bit narrower to fit the new content
#  - hyperlinks are rebuilt so rId1..rId7 line up with the 7 remaining rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$NEW_TS = "2025-11-05 06:27:26"

# --- 1. Drop the obsolete rows (9-19) entirely, shifting nothing else ---
$ws.Range("A9:H19").EntireRow.Delete()

# --- 2. Refresh the capture timestamp for every remaining data row (2-8) ---
$ws.Range("A2:A8").Value = $NEW_TS

# --- 3. Rows 2-4 keep their content (only the timestamp changed above). ---

# --- 4. Row 5: new listing ---
$ws.Range("B5").Value = "<Next.js、バックエンド開発> ガントチャートアプリの改修製造"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5427011"
$ws.Range("G5").Value = 225
$ws.Range("H5").Value = "🔥Next.js ◆開発 ◇アプリ"

# --- 5. Row 6: new listing ---
$ws.Range("B6").Value = "デフォルトカメラ機能を活用したアプリ開発"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5427397"
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = "◆開発 ◇アプリ"

# --- 6. Row 7: new listing ---
$ws.Range("B7").Value = "弥生販売 得意先台帳登録 商品登録 売上伝票作成ツールのご相談"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5427338"
$ws.Range("G7").Value = 73
$ws.Range("H7").Value = "◆ツール"

# --- 7. Row 8: new listing (no skill-summary column this time) ---
$ws.Range("B8").Value = "【継続依頼あり】GASやn8nのオンラインセミナー研修講師を募集!"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5427459"
$ws.Range("G8").Value = 13
$ws.Range("H8").ClearContents()

# --- 8. Rebuild the hyperlinks collection so it only covers F2:F8 ---
#     (the per-item Hyperlink.Delete() is a no-op in this engine, so the
#      collection has to be cleared and re-populated wholesale)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5416301")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5420440")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5416328")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5427011")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5427397")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5427338")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5427459")

# Hyperlinks.Add() re-styles the target cell with a brand-new cellXfs entry;
# fold it back onto the original shared "Hyperlink" style so F2:F8 keep the
# same style index (s="1") they had before.
$ws.Range("F2:F8").Style = "Hyperlink"

# --- 9. Column widths: D (4) 32->28, H (8) 27->19 ---
# ColumnWidth uses Excel's padded character-width units, not the raw OOXML
# "width" value, so aim at the input that rounds back to the exact target.
$ws.Columns.Item(4).ColumnWidth = 27.17
$ws.Columns.Item(8).ColumnWidth = 18.17
